# Updated symbol list (price / 1h volume change) refresh for cryptos sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value. Values are written as text
# (leading apostrophe forces text interpretation so "317.33" stays a
# string and "-3.45%" isn't reinterpreted as a percentage number), then
# the style is reset to "Normal" so no extra number-format/style is left
# on the cell (matches the original inline-string cells, which carry no
# special style).
$updates = [ordered]@{
    "D2"  = "317.33";    "E2"  = "-3.45%"
    "E3"  = "-5.03%"
    "D4"  = "5.195";     "E4"  = "-3.24%"
    "D5"  = "0.08104";   "E5"  = "-3.35%"
    "D6"  = "4.371";     "E6"  = "-1.51%"
    "D7"  = "1.746";     "E7"  = "-10.34%"
    "D8"  = "0.9283";    "E8"  = "-4.86%"
    "D9"  = "0.1121";    "E9"  = "0.34%"
    "D10" = "0.1855";    "E10" = "-2.41%"
    "D11" = "0.09268";   "E11" = "-3.99%"
    "D12" = "0.04586";   "E12" = "-0.54%"
    "D13" = "7.389";     "E13" = "-19.34%"
    "D14" = "0.1052";    "E14" = "-1.07%"
    "D15" = "0.001280";  "E15" = "-0.87%"
    "D16" = "0.005999";  "E16" = "-1.88%"
    "D17" = "3.344";     "E17" = "-1.78%"
    "E18" = "1.47%"
    "D19" = "0.3389";    "E19" = "1.76%"
    "D20" = "0.1385";    "E20" = "1.04%"
    "D21" = "0.2605";    "E21" = "2.15%"
    "D22" = "0.04182";   "E22" = "0.80%"
    "D23" = "0.001244";  "E23" = "-4.01%"
    "D24" = "0.004243";  "E24" = "-3.62%"
    "D25" = "0.0001224"; "E25" = "-5.97%"
    "D26" = "0.0002988"; "E26" = "0.04%"
    "D38" = "0.02582";   "E38" = "-2.79%"
    "D39" = "0.05469";   "E39" = "-2.77%"
    "D40" = "0.008031";  "E40" = "2.34%"
    "D41" = "0.1388";    "E41" = "-1.77%"
    "D42" = "0.006543";  "E42" = "-11.10%"
    "D43" = "0.002086";  "E43" = "-1.25%"
    "D44" = "0.008233";  "E44" = "4.11%"
    "D45" = "0.3448";    "E45" = "-1.95%"
    "D46" = "0.00006749";    "E46" = "-2.10%"
    "D47" = "0.00000000752"; "E47" = "0.19%"
    "D48" = "0.003397";  "E48" = "-3.34%"
    "D49" = "0.004113"
    "D50" = "0.00002107"; "E50" = "0.19%"
    "D51" = "0.0002006";  "E51" = "0.19%"
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.Value = "'" + $updates[$cellRef]
    $cell.Style = "Normal"
}
